$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = "26.662.17"
$ws.Range("E2").Value = "  +1.55%  "

# Row 3
$ws.Range("D3").Value = "1.634.75"
$ws.Range("E3").Value = "  +1.06%  "

# Row 4
$ws.Range("E4").Value = "  -0.06%  "

# Row 5
$ws.Range("D5").Value = "'213.44"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +0.70%  "

# Row 6
$ws.Range("D6").Value = "'0.506"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +4.15%  "

# Row 7
$ws.Range("E7").Value = "  +0.02%  "

# Row 8
$ws.Range("E8").Value = "  +2.65%  "

# Row 9
$ws.Range("E9").Value = "  +1.45%  "

# Row 10
$ws.Range("D10").Value = "'19.28"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +2.85%  "

# Row 11
$ws.Range("D11").Value = "'0.0844"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +3.48%  "

# Row 12
$ws.Range("D12").Value = "1.864.05"
$ws.Range("E12").Value = "  +1.15%  "

# Row 13
$ws.Range("D13").Value = "1.675.25"
$ws.Range("E13").Value = "  +3.52%  "

# Row 14
$ws.Range("E14").Value = "  +2.82%  "

# Row 15
$ws.Range("D15").Value = "'0.527"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +2.00%  "

# Row 16
$ws.Range("D16").Value = "26.671.25"
$ws.Range("E16").Value = "  +1.55%  "

# Row 17
$ws.Range("D17").Value = "'63.54"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +2.12%  "

# Row 18
$ws.Range("D18").Value = "0.0₃0744"
$ws.Range("E18").Value = "  +2.54%  "

# Row 19
$ws.Range("D19").Value = "'219.34"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +9.23%  "

# Row 20
$ws.Range("E20").Value = "  -0.01%  "

# Row 21
$ws.Range("D21").Value = "'4.30"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +0.84%  "

# Row 22
$ws.Range("B22").Value = "Chainlink"
$ws.Range("C22").Value = "https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link"
$ws.Range("D22").Value = "'6.20"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +2.81%  "

# Row 23
$ws.Range("B23").Value = "Avalanche"
$ws.Range("C23").Value = "https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax"
$ws.Range("D23").Value = "'9.43"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +1.37%  "

# Row 24
$ws.Range("D24").Value = "'1.93"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +1.62%  "

# Row 25
$ws.Range("D25").Value = "'148.11"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +2.77%  "

# Row 26
$ws.Range("E26").Value = "  +0.01%  "

# Row 27
$ws.Range("E27").Value = "  +1.53%  "

# Row 28
$ws.Range("E28").Value = "  +6.10%  "

# Row 29
$ws.Range("E29").Value = "  +2.29%  "

# Row 30
$ws.Range("D30").Value = "'0.0510"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -0.80%  "

# Row 31
$ws.Range("E31").Value = "  +0.02%  "

# Row 32
$ws.Range("E32").Value = "  +4.81%  "

# Row 33
$ws.Range("E33").Value = "  +2.83%  "

# Row 34
$ws.Range("E34").Value = "  +1.85%  "

# Row 35
$ws.Range("E35").Value = "  -0.28%  "

# Row 36
$ws.Range("D36").Value = "1.211.73"
$ws.Range("E36").Value = "  +3.20%  "

# Row 37
$ws.Range("E37").Value = "  +5.64%  "

# Row 38
$ws.Range("D38").Value = "'0.810"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +1.13%  "

# Row 39
$ws.Range("E39").Value = "  +0.06%  "

# Row 40
$ws.Range("D40").Value = "'0.505"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +2.35%  "

# Row 41
$ws.Range("E41").Value = "  -1.07%  "

# Row 42
$ws.Range("D42").Value = "'5.43"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +1.85%  "

# Row 43
$ws.Range("D43").Value = "'0.793"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +0.27%  "

# Row 44
$ws.Range("D44").Value = "1.773.59"
$ws.Range("E44").Value = "  +1.05%  "

# Row 45
$ws.Range("D45").Value = "'93.21"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +0.72%  "

# Row 46
$ws.Range("E46").Value = "  +2.00%  "

# Row 47
$ws.Range("D47").Value = "'54.86"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +2.39%  "

# Row 48
$ws.Range("D48").Value = "'0.0513"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +1.00%  "

# Row 49
$ws.Range("D49").Value = "'7.71"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +6.77%  "

# Row 50
$ws.Range("E50").Value = "  +0.33%  "

# Row 51
$ws.Range("E51").Value = "  +0.27%  "
